# Generate Report for Handback
# For the "a3d662e1-c28a-40a3-aee1-6bf92006a07a" file row (row 8) on both the
# zh-cn and de-de sheets, the handback pass found a stale handback package:
# it fills in the "Latest Target File" hyperlink, "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns, and widens the
# Error Detail column so the message is readable.

$wb = $excel.ActiveWorkbook

$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f971c5c2a89bce4797d4d5bfa90d7e93645909/e2e/a3d662e1-c28a-40a3-aee1-6bf92006a07a.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7521fa76a89832be0ae3c4a4b273cf6d8218310f/e2e/a3d662e1-c28a-40a3-aee1-6bf92006a07a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f971c5c2a89bce4797d4d5bfa90d7e93645909/e2e/a3d662e1-c28a-40a3-aee1-6bf92006a07a.md."

# --- zh-cn sheet ---------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $latestHandbackUrl, "", "", "a3d662e1-c28a-40a3-aee1-6bf92006a07a.md")
$wsZh.Range("J8").Value = "a3d662e1-c28a-40a3-aee1-6bf92006a07a.d9603b810c752896d36a88427914467b80e91aa2.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-09-07 02:53:42"
$wsZh.Range("P8").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet -----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $latestHandbackUrl, "", "", "a3d662e1-c28a-40a3-aee1-6bf92006a07a.md")
$wsDe.Range("J8").Value = "a3d662e1-c28a-40a3-aee1-6bf92006a07a.d9603b810c752896d36a88427914467b80e91aa2.de-de.xlf"
$wsDe.Range("K8").Value = "2016-09-07 02:53:50"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.16666666666667
